$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price column D, Volume(1h) column E)
# Values are written as text (apostrophe-prefixed) to match the
# source sheet convention, where Price/Volume are stored as literal
# strings rather than parsed numbers.
$ws.Range("D2").Value = "'294.58"
$ws.Range("E2").Value = "'-3.77%"
$ws.Range("D3").Value = "'31.27"
$ws.Range("E3").Value = "'-2.59%"
$ws.Range("E4").Value = "'-3.45%"
$ws.Range("D5").Value = "'0.07367"
$ws.Range("E5").Value = "'-0.02%"
$ws.Range("D6").Value = "'7.683"
$ws.Range("E6").Value = "'-2.01%"
$ws.Range("D7").Value = "'3.757"
$ws.Range("E7").Value = "'-1.02%"
$ws.Range("D8").Value = "'1.632"
$ws.Range("E8").Value = "'9.72%"
$ws.Range("D9").Value = "'0.9240"
$ws.Range("E9").Value = "'1.48%"
$ws.Range("D10").Value = "'0.1669"
$ws.Range("E10").Value = "'-1.23%"
$ws.Range("D11").Value = "'0.07178"
$ws.Range("E11").Value = "'-4.81%"
$ws.Range("D12").Value = "'0.07944"
$ws.Range("E12").Value = "'-1.14%"
$ws.Range("D13").Value = "'0.02995"
$ws.Range("E13").Value = "'-1.09%"
$ws.Range("D14").Value = "'0.09892"
$ws.Range("E14").Value = "'-1.05%"
$ws.Range("D15").Value = "'0.001491"
$ws.Range("E15").Value = "'-0.99%"
$ws.Range("D16").Value = "'0.006209"
$ws.Range("E16").Value = "'0.28%"
$ws.Range("D17").Value = "'3.455"
$ws.Range("E17").Value = "'-0.78%"
$ws.Range("E18").Value = "'-0.17%"
$ws.Range("E19").Value = "'-0.15%"
$ws.Range("D20").Value = "'0.1335"
$ws.Range("E20").Value = "'-0.80%"
$ws.Range("D21").Value = "'4.556"
$ws.Range("E21").Value = "'5.07%"
$ws.Range("D22").Value = "'0.04630"
$ws.Range("E22").Value = "'1.38%"
$ws.Range("D23").Value = "'0.1552"
$ws.Range("E23").Value = "'-5.81%"
$ws.Range("E24").Value = "'-1.04%"
$ws.Range("E25").Value = "'-0.49%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'-1.07%"
$ws.Range("E27").Value = "'6.83%"
$ws.Range("D39").Value = "'0.01685"
$ws.Range("E39").Value = "'1.90%"
$ws.Range("D40").Value = "'0.04404"
$ws.Range("E40").Value = "'-2.24%"
$ws.Range("D41").Value = "'0.007124"
$ws.Range("E41").Value = "'-2.46%"
$ws.Range("E42").Value = "'-1.87%"
$ws.Range("D43").Value = "'0.002100"
$ws.Range("E43").Value = "'-8.08%"
$ws.Range("D44").Value = "'0.01102"
$ws.Range("E44").Value = "'-22.39%"
$ws.Range("D45").Value = "'0.00006004"
$ws.Range("E45").Value = "'-1.31%"
$ws.Range("D46").Value = "'1.918"
$ws.Range("E46").Value = "'1.34%"
$ws.Range("E47").Value = "'-16.09%"
